$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) columns with the latest crypto snapshot values.
# D-column values that look like plain numbers are briefly forced to Text so Excel
# keeps the original fixed-precision string (e.g. "333.61") instead of normalizing it,
# then the number format is restored to General to match the source formatting.

$ws.Range("D2").Value = "26.699.15"
$ws.Range("E2").Value = "  +7.26%  "

$ws.Range("D3").Value = "1.743.07"
$ws.Range("E3").Value = "  +5.01%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.61"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +6.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.33%  "

$ws.Range("E7").Value = "  +3.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.40"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +4.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3405"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +4.90%  "

$ws.Range("E10").Value = "  +6.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07483"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +6.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9981"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.474"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +7.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.45"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +5.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.116"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +8.40%  "

$ws.Range("D16").Value = "1.743.55"
$ws.Range("E16").Value = "  +5.15%  "

$ws.Range("E17").Value = "  +4.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06693"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.71"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +6.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9990"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.78"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +7.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.192"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +5.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.10"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +4.48%  "

$ws.Range("D24").Value = "26.702.79"
$ws.Range("E24").Value = "  +7.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.438"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.478"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +2.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.425"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +19.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.41"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +4.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.69"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +6.24%  "

$ws.Range("D30").Value = "1.931.04"
$ws.Range("E30").Value = "  +4.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.36"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +6.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.128"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.91%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.123"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +6.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08652"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.702"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.04"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +6.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.458"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +6.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02367"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +5.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06321"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +4.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2182"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +6.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.609"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +4.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.234"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -3.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6260"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +6.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.29"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +12.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.926"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +4.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6064"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +8.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.38"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +3.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.067"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +7.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07295"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +4.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.88"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +4.78%  "
